$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 3.558321333333333
$ws.Cells.Item(2, 8).Value = 10.674964
$ws.Cells.Item(2, 9).Value = 0.3039644761000113
$ws.Cells.Item(2, 10).Value = 0.3039644761000113
$ws.Cells.Item(2, 13).Value = 0.5569716666666666
$ws.Cells.Item(2, 14).Value = 1.670915
$ws.Cells.Item(2, 15).Value = 0.3425729542218473
$ws.Cells.Item(2, 16).Value = 0.3425729542218473
$ws.Cells.Item(2, 17).Value = 1.981884163562222
$ws.Cells.Item(2, 18).Value = 17.83695747206
$ws.Cells.Item(2, 19).Value = 0.104130008556077
$ws.Cells.Item(2, 20).Value = 0.104130008556077

$ws.Cells.Item(3, 7).Value = 3.558321333333333
$ws.Cells.Item(3, 8).Value = 10.674964
$ws.Cells.Item(3, 9).Value = 0.3039644761000113
$ws.Cells.Item(3, 10).Value = 0.3039644761000113
$ws.Cells.Item(3, 13).Value = 0.6936943333333333
$ws.Cells.Item(3, 14).Value = 2.081083
$ws.Cells.Item(3, 15).Value = 0.4266660789393025
$ws.Cells.Item(3, 16).Value = 0.4266660789393025
$ws.Cells.Item(3, 17).Value = 2.468387345112444
$ws.Cells.Item(3, 18).Value = 22.215486106012
$ws.Cells.Item(3, 19).Value = 0.1296913311544312
$ws.Cells.Item(3, 20).Value = 0.1296913311544311

$ws.Cells.Item(4, 7).Value = 3.558321333333333
$ws.Cells.Item(4, 8).Value = 10.674964
$ws.Cells.Item(4, 9).Value = 0.3039644761000113
$ws.Cells.Item(4, 10).Value = 0.3039644761000113
$ws.Cells.Item(4, 13).Value = 0.3751823333333333
$ws.Cells.Item(4, 14).Value = 1.125547
$ws.Cells.Item(4, 15).Value = 0.2307609668388503
$ws.Cells.Item(4, 16).Value = 0.2307609668388503
$ws.Cells.Item(4, 17).Value = 1.335019300589778
$ws.Cells.Item(4, 18).Value = 12.015173705308
$ws.Cells.Item(4, 19).Value = 0.07014313638950322
$ws.Cells.Item(4, 20).Value = 0.0701431363895032

$ws.Cells.Item(5, 7).Value = 5.383140666666667
$ws.Cells.Item(5, 8).Value = 16.149422
$ws.Cells.Item(5, 9).Value = 0.4598470400038817
$ws.Cells.Item(5, 10).Value = 0.4598470400038817
$ws.Cells.Item(5, 13).Value = 0.5569716666666666
$ws.Cells.Item(5, 14).Value = 1.670915
$ws.Cells.Item(5, 15).Value = 0.3425729542218473
$ws.Cells.Item(5, 16).Value = 0.3425729542218473
$ws.Cells.Item(5, 17).Value = 2.998256829014444
$ws.Cells.Item(5, 18).Value = 26.98431146113
$ws.Cells.Item(5, 19).Value = 0.1575311589843018
$ws.Cells.Item(5, 20).Value = 0.1575311589843017

$ws.Cells.Item(6, 7).Value = 5.383140666666667
$ws.Cells.Item(6, 8).Value = 16.149422
$ws.Cells.Item(6, 9).Value = 0.4598470400038817
$ws.Cells.Item(6, 10).Value = 0.4598470400038817
$ws.Cells.Item(6, 13).Value = 0.6936943333333333
$ws.Cells.Item(6, 14).Value = 2.081083
$ws.Cells.Item(6, 15).Value = 0.4266660789393025
$ws.Cells.Item(6, 16).Value = 0.4266660789393025
$ws.Cells.Item(6, 17).Value = 3.734254176002889
$ws.Cells.Item(6, 18).Value = 33.608287584026
$ws.Cells.Item(6, 19).Value = 0.1962011334703008
$ws.Cells.Item(6, 20).Value = 0.1962011334703007

$ws.Cells.Item(7, 7).Value = 5.383140666666667
$ws.Cells.Item(7, 8).Value = 16.149422
$ws.Cells.Item(7, 9).Value = 0.4598470400038817
$ws.Cells.Item(7, 10).Value = 0.4598470400038817
$ws.Cells.Item(7, 13).Value = 0.3751823333333333
$ws.Cells.Item(7, 14).Value = 1.125547
$ws.Cells.Item(7, 15).Value = 0.2307609668388503
$ws.Cells.Item(7, 16).Value = 0.2307609668388503
$ws.Cells.Item(7, 17).Value = 2.019659275981556
$ws.Cells.Item(7, 18).Value = 18.176933483834
$ws.Cells.Item(7, 19).Value = 0.1061147475492792
$ws.Cells.Item(7, 20).Value = 0.1061147475492792

$ws.Cells.Item(8, 7).Value = 2.764910333333333
$ws.Cells.Item(8, 8).Value = 8.294730999999999
$ws.Cells.Item(8, 9).Value = 0.2361884838961071
$ws.Cells.Item(8, 10).Value = 0.236188483896107
$ws.Cells.Item(8, 13).Value = 0.5569716666666666
$ws.Cells.Item(8, 14).Value = 1.670915
$ws.Cells.Item(8, 15).Value = 0.3425729542218473
$ws.Cells.Item(8, 16).Value = 0.3425729542218473
$ws.Cells.Item(8, 17).Value = 1.539976716540555
$ws.Cells.Item(8, 18).Value = 13.859790448865
$ws.Cells.Item(8, 19).Value = 0.0809117866814686
$ws.Cells.Item(8, 20).Value = 0.08091178668146859

$ws.Cells.Item(9, 7).Value = 2.764910333333333
$ws.Cells.Item(9, 8).Value = 8.294730999999999
$ws.Cells.Item(9, 9).Value = 0.2361884838961071
$ws.Cells.Item(9, 10).Value = 0.236188483896107
$ws.Cells.Item(9, 13).Value = 0.6936943333333333
$ws.Cells.Item(9, 14).Value = 2.081083
$ws.Cells.Item(9, 15).Value = 0.4266660789393025
$ws.Cells.Item(9, 16).Value = 0.4266660789393025
$ws.Cells.Item(9, 17).Value = 1.918002630408111
$ws.Cells.Item(9, 18).Value = 17.262023673673
$ws.Cells.Item(9, 19).Value = 0.1007736143145706
$ws.Cells.Item(9, 20).Value = 0.1007736143145706

$ws.Cells.Item(10, 7).Value = 2.764910333333333
$ws.Cells.Item(10, 8).Value = 8.294730999999999
$ws.Cells.Item(10, 9).Value = 0.2361884838961071
$ws.Cells.Item(10, 10).Value = 0.236188483896107
$ws.Cells.Item(10, 13).Value = 0.3751823333333333
$ws.Cells.Item(10, 14).Value = 1.125547
$ws.Cells.Item(10, 15).Value = 0.2307609668388503
$ws.Cells.Item(10, 16).Value = 0.2307609668388503
$ws.Cells.Item(10, 17).Value = 1.037345510317444
$ws.Cells.Item(10, 18).Value = 9.336109592856999
$ws.Cells.Item(10, 19).Value = 0.05450308290006789
$ws.Cells.Item(10, 20).Value = 0.05450308290006788
